$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("11-4-22")
$ws.Activate()

# Row 4 - Designing prototype -> Entity, Date model and relationships / Training Head and Co-ordinator Dessign flow
$ws.Range("B4").Value = "Entity, Date model and relationships"
$ws.Range("C4").Value = "Training Head and Co-ordinator Dessign flow "

# Row 8 - update B8, C8, E8, F8
$ws.Range("B8").Value = "Entity data model operation for Trainer"
$ws.Range("C8").Value = "Entity datamodel ,artributes,relationship for trainer"
$ws.Range("E8").Value = "Discussion with team - 20 mins, Worked on entity data model for trainer -3 hours ,meeting with Rafi- 100 mins"
$ws.Range("F8").Value = "others(lunch & tea break)-1.30hour                       "

# Row 9 - was ABSENT for all cells, now filled in (D9 stays blank)
$ws.Range("B9").Value = "Entity data model "
$ws.Range("C9").Value = "entered sample data for head , co ordinator, trainee , trainer."
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = "discussion with team -20min, reviewed the head and co ordinator UI, meeting rafi 100 mins , entering sample data for head,co ordinator, trainer , trainee - 3 hrs, editing the changes in UI ( head , co ordinator, trainer, trainee)."
$ws.Range("F9").Value = "lunch and break 1.30hrs"

# Row 11 - B11
$ws.Range("B11").Value = "Entity Data model and its Relationship"

# Update the view state: scrolled to row 6, selection on F9
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F9").Select()
